$d = $word.ActiveDocument

# The scraped agenda page text, one array element per visual line of the
# source PDF/HTML capture. Joining with a manual line break (vertical-tab,
# i.e. Shift+Enter / <w:br/>) keeps everything inside a single paragraph,
# matching how the page was laid out.
$lines = @(
    "cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4415",
    "6/23/23, 5:41 PM",
    "ORLEANS",
    "NEW O",
    "CITY OF",
    "*",
    "LOUISIANA",
    "CRIMINAL JUSTICE COMMITTEE",
    "MEETING AGENDA",
    "TUESDAY, FEBRUARY 7, 2023",
    "09:30 A.M.",
    "Quarterly report: Update on public safety initiatives and legislative, financial or",
    "3.",
    "operational issues that have been identified:",
    "Orleans Parish Criminal District Court, Judge Robin Pittman and/or representative",
    "A.",
    "1. Roll Call",
    "2. Approval of the minutes from the January 23, 2023 meeting",
    "E. Juvenile Justice Intervention Center, Dichelle Williams and/or representative",
    "TJJIC Presentation",
    "B. Orleans Parish Juvenile Court, Judge Ranord Darensburg and/or representative",
    "Orleans Parish Juvenile Court Presentation",
    "Orleans Parish District Attorney, District Attorney Jason Williams and/or",
    "C.",
    "representative",
    "D.A.'s Presentation",
    "D. New Orleans Police Department, Superintendent Michelle Woodfork and/or",
    "representative",
    "NOPD Presentation",
    "Office of Independent Police Monitor, Stella Cziment and/or representative",
    "F.",
    "OIPM Presentation",
    "G. Orleans Parish Sheriff, Sheriff Susan Hutson and/or representative",
    "https://cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4415",
    "1/2",
    "cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4415",
    "6/23/23, 5:41 PM",
    "Sheriff's Presentation",
    "4. Adjournment",
    "Public Comment",
    "2/2",
    "https://cityofno.granicus.com/GeneratedAgendaViewer.php?view_id=42&clip_id=4415"
)

$lineBreak = [char]11   # manual line break -> <w:br/>
$paraBreak = [char]13   # paragraph mark -> new <w:p>

$bodyText = [string]::Join($lineBreak, $lines)
$footer = $lineBreak + "---------- End of Page 1 ----------" + $lineBreak

$d.Paragraphs(1).Range.Text = $bodyText + $paraBreak + $footer
